# horarios-141-completo.xlsx update: 30/12 17:56 LP1912+6203+6173
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: LP1912  (columns: A=(meta/blank) B=Hora_Scrap C=Hora_Llegada
#          D=Linea E=Minutos(num) F=Parada G=Fecha)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 14:56:30"
$ws1.Cells.Item(3,1).Value = "Total filas: 319"

$sheet1Rows = @(
    @("14:56:19","15:05","10_OLMOS",9,"LP1912","30/12/2025"),
    @("14:56:19","15:05","23_HERNANDEZ",9,"LP1912","30/12/2025"),
    @("14:56:19","15:07","16_SANTA ANA",11,"LP1912","30/12/2025"),
    @("14:56:19","15:20","15_ABASTO",24,"LP1912","30/12/2025"),
    @("14:56:19","15:21","26_HERNANDEZ",25,"LP1912","30/12/2025"),
    @("14:56:19","15:27","16_SANTA ANA",31,"LP1912","30/12/2025"),
    @("14:56:19","15:32","84_COLONIA URQUIZA-ESC 49",36,"LP1912","30/12/2025"),
    @("14:56:19","15:42","10_OLMOS",46,"LP1912","30/12/2025"),
    @("14:56:19","15:46","14_ABASTO",50,"LP1912","30/12/2025"),
    @("14:56:19","15:54","11_ETCHEVERRY",58,"LP1912","30/12/2025"),
    @("14:56:19","15:54","23_HERNANDEZ",58,"LP1912","30/12/2025"),
    @("14:56:19","16:01","10_OLMOS",65,"LP1912","30/12/2025"),
    @("14:56:19","16:07","23_HERNANDEZ",71,"LP1912","30/12/2025"),
    @("14:56:19","16:20","215C_EL PATO",84,"LP1912","30/12/2025"),
    @("14:56:19","16:21","26_HERNANDEZ",85,"LP1912","30/12/2025"),
    @("14:56:19","16:26","14_ABASTO",90,"LP1912","30/12/2025"),
    @("14:56:19","16:32","11_ETCHEVERRY",96,"LP1912","30/12/2025")
)

$r = 304
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r,2).Value = $row[0]
    $ws1.Cells.Item($r,3).Value = $row[1]
    $ws1.Cells.Item($r,4).Value = $row[2]
    $ws1.Cells.Item($r,5).Value = $row[3]
    $ws1.Cells.Item($r,6).Value = $row[4]
    $ws1.Cells.Item($r,7).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet 2: LP1912-215  (columns: A=(meta/blank) B=Fecha C=Hora_Scrap
#          D=Hora_Llegada E=Linea F=Minutos(num) G=Parada)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 14:56:30"
$ws2.Cells.Item(3,1).Value = "Total filas: 24"

$ws2.Cells.Item(25,2).Value = "30/12/2025"
$ws2.Cells.Item(25,3).Value = "14:56:19"
$ws2.Cells.Item(25,4).Value = "16:20"
$ws2.Cells.Item(25,5).Value = "215C_EL PATO"
$ws2.Cells.Item(25,6).Value = 84
$ws2.Cells.Item(25,7).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 3: 6203-6173  (columns: A=(meta/blank) B=Fecha C=Hora_Scrap
#          D=Hora_Llegada E=Linea F=Minutos(num) G=Parada)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 14:56:30"
$ws3.Cells.Item(3,1).Value = "Total filas: 45"

$ws3.Cells.Item(46,2).Value = "30/12/2025"
$ws3.Cells.Item(46,3).Value = "14:56:30"
$ws3.Cells.Item(46,4).Value = "15:34"
$ws3.Cells.Item(46,5).Value = "215A_LA PLATA"
$ws3.Cells.Item(46,6).Value = 38
$ws3.Cells.Item(46,7).Value = "L6173"
